$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A1").Value = 2253.0246659930554
$ws.Range("B1").Value = 1484.5375104173265
$ws.Range("C1").Value = 1458.5736134654164
$ws.Range("A2").Value = 2245.2662449656004
$ws.Range("B2").Value = 1471.2404910463865
$ws.Range("C2").Value = 1380.7460134102278
$ws.Range("A3").Value = 2336.1676700863045
$ws.Range("B3").Value = 1617.0565113740436
$ws.Range("C3").Value = 1523.1838345469969
$ws.Range("A4").Value = 2342.9978193578449
$ws.Range("B4").Value = 1761.5124566489699
$ws.Range("C4").Value = 1689.8979045080016
$ws.Range("A5").Value = 2434.1217742752733
$ws.Range("B5").Value = 1672.2818283659492
$ws.Range("C5").Value = 1684.3537791446468
$ws.Range("A6").Value = 2391.6738806367443
$ws.Range("B6").Value = 1801.4606032007052
$ws.Range("C6").Value = 1844.3669383967613
$ws.Range("A7").Value = 2137.565899469203
$ws.Range("B7").Value = 1574.3978801186047
$ws.Range("C7").Value = 1495.5658458013481
$ws.Range("A8").Value = 2229.037756518248
$ws.Range("B8").Value = 1669.9262810187711
$ws.Range("C8").Value = 1623.2863319398716
$ws.Range("A9").Value = 2482.1508319710756
$ws.Range("B9").Value = 1785.8529223760083
$ws.Range("C9").Value = 1532.2339656633517
$ws.Range("A10").Value = 2137.6608258863239
$ws.Range("B10").Value = 1367.6007277853821
$ws.Range("C10").Value = 1355.2244003357537
$ws.Range("A11").Value = 1916.785965067352
$ws.Range("B11").Value = 1487.2775584793419
$ws.Range("C11").Value = 1340.4566603589826
$ws.Range("A12").Value = 2624.0696799211869
$ws.Range("B12").Value = 2214.4940293367436
$ws.Range("C12").Value = 2007.622886469951
$ws.Range("A13").Value = 2434.1640767591011
$ws.Range("B13").Value = 1772.1942831992394
$ws.Range("C13").Value = 1801.0800102327964
$ws.Range("A14").Value = 2533.7209044577908
$ws.Range("B14").Value = 1892.4211717170908
$ws.Range("C14").Value = 1672.6306134579916
$ws.Range("A15").Value = 2608.6793612619822
$ws.Range("B15").Value = 2053.6014331011174
$ws.Range("C15").Value = 1940.6135712709113
$ws.Range("A16").Value = 2247.0467798670688
$ws.Range("B16").Value = 1531.9959506533103
$ws.Range("C16").Value = 1285.6398882655087
